# Backup before import design and part overhaul
# Applies: fix H3 value, and append two new data rows (rows 4 and 5)
# mirroring the structure of existing rows, with new Destination Well
# values "A3" and "A4".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the Transfer Volume for row 3 (UID 2)
$ws.Range("H3").Value = 2625

# Repeating column literals shared by every data row.
$colB = "level 1 6RES source plate"
$colC = "6RES_AQ_BP"
$colD = "A1"
$colE = "384-Well Level 1 MoClo output plate"
$colF = "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)"
$colI = "Deionised water"

# New row 4 -> UID 3, Destination Well "A3"
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = $colB
$ws.Range("C4").Value = $colC
$ws.Range("D4").Value = $colD
$ws.Range("E4").Value = $colE
$ws.Range("F4").Value = $colF
$ws.Range("G4").Value = "A3"
$ws.Range("H4").Value = 2875
$ws.Range("I4").Value = $colI

# New row 5 -> UID 4, Destination Well "A4"
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = $colB
$ws.Range("C5").Value = $colC
$ws.Range("D5").Value = $colD
$ws.Range("E5").Value = $colE
$ws.Range("F5").Value = $colF
$ws.Range("G5").Value = "A4"
$ws.Range("H5").Value = 2625
$ws.Range("I5").Value = $colI
